$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 111.61539
$ws.Range("I33").Value = 104.5
$ws.Range("J33").Value = 123
$ws.Range("K33").Value = 104.5
$ws.Range("L33").Value = 123
$ws.Range("M33").Value = 124.5
$ws.Range("N33").Value = -581

$ws.Range("H132").Value = 644948.2
$ws.Range("I132").Value = 810020.7
$ws.Range("J132").Value = 25926.5
$ws.Range("K132").Value = 2430062.1
$ws.Range("L132").Value = 77779.5
$ws.Range("M132").Value = -2427532.1
$ws.Range("N132").Value = -82839.5

$ws.Range("H137").Value = 1309.1364
$ws.Range("I137").Value = 988.06665
$ws.Range("J137").Value = 1997.1428
$ws.Range("K137").Value = 2964.19995
$ws.Range("L137").Value = 5991.428400000001
$ws.Range("M137").Value = -414.1999500000002
$ws.Range("N137").Value = -11091.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 64244.312
$ws.Range("I2").Value = 101930.9
$ws.Range("J2").Value = 1433.3334
$ws.Range("K2").Value = 101930.9
$ws.Range("L2").Value = 1433.3334
$ws.Range("M2").Value = -101817.9
$ws.Range("N2").Value = -1659.3334

$ws.Range("H64").Value = 15000
$ws.Range("I64").Value = 15000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 15000
$ws.Range("L64").Value = $null
$ws.Range("M64").Value = -14752
$ws.Range("N64").Value = $null

$ws.Range("H67").Value = 15000
$ws.Range("I67").Value = 15000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = $null
$ws.Range("M67").Value = -14142
$ws.Range("N67").Value = $null

$ws.Range("H74").Value = 4837.3335
$ws.Range("I74").Value = 879.75
$ws.Range("J74").Value = 26999.8
$ws.Range("K74").Value = 879.75
$ws.Range("L74").Value = 26999.8
$ws.Range("M74").Value = -5.75
$ws.Range("N74").Value = -28747.8

$ws.Range("H77").Value = 4837.3335
$ws.Range("I77").Value = 879.75
$ws.Range("J77").Value = 26999.8
$ws.Range("K77").Value = 4398.75
$ws.Range("L77").Value = 134999
$ws.Range("M77").Value = -30.75
$ws.Range("N77").Value = -143735

$ws.Range("H92").Value = 30037
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 30037
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 30037
$ws.Range("N92").Value = -35029

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = $null
$ws.Range("N93").Value = $null

$ws.Range("H94").Value = 24155
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 24155
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 24155
$ws.Range("N94").Value = -25957

$ws.Range("H116").Value = 64244.312
$ws.Range("I116").Value = 101930.9
$ws.Range("J116").Value = 1433.3334
$ws.Range("K116").Value = 101930.9
$ws.Range("L116").Value = 1433.3334
$ws.Range("M116").Value = -99636.89999999999
$ws.Range("N116").Value = -6021.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 64244.312
$ws.Range("I3").Value = 101930.9
$ws.Range("J3").Value = 1433.3334
$ws.Range("K3").Value = 101930.9
$ws.Range("L3").Value = 1433.3334
$ws.Range("M3").Value = -101816.9
$ws.Range("N3").Value = -1661.3334

$ws.Range("H134").Value = 4022.4583
$ws.Range("I134").Value = 2920.8667
$ws.Range("J134").Value = 5858.4443
$ws.Range("K134").Value = 8762.6001
$ws.Range("L134").Value = 17575.3329
$ws.Range("M134").Value = -6227.6001
$ws.Range("N134").Value = -22645.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -1574

$ws.Range("H31").Value = 2189.35
$ws.Range("I31").Value = 1355.591
$ws.Range("J31").Value = 3208.389
$ws.Range("K31").Value = 1355.591
$ws.Range("L31").Value = 3208.389
$ws.Range("M31").Value = -1060.591
$ws.Range("N31").Value = -3798.389

$ws.Range("H34").Value = 2189.35
$ws.Range("I34").Value = 1355.591
$ws.Range("J34").Value = 3208.389
$ws.Range("K34").Value = 1355.591
$ws.Range("L34").Value = 3208.389
$ws.Range("M34").Value = -1153.591
$ws.Range("N34").Value = -3612.389

$ws.Range("H107").Value = 791.2143
$ws.Range("I107").Value = 422.8
$ws.Range("J107").Value = 995.8889
$ws.Range("K107").Value = 422.8
$ws.Range("L107").Value = 995.8889
$ws.Range("M107").Value = 1497.2
$ws.Range("N107").Value = -4835.8889

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -5340

$ws.Range("H134").Value = 4844.6875
$ws.Range("I134").Value = 2419
$ws.Range("J134").Value = 6300.1
$ws.Range("K134").Value = 7257
$ws.Range("L134").Value = 18900.3
$ws.Range("M134").Value = -4722
$ws.Range("N134").Value = -23970.3

$ws.Range("H135").Value = 42719.168
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 42719.168
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 42719.168
$ws.Range("N135").Value = -52859.168

$ws.Range("H141").Value = 1356750
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 1356750
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 1356750
$ws.Range("N141").Value = -1367110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 25625.25
$ws.Range("I68").Value = 25625.25
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 76875.75
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = $null

$ws.Range("H71").Value = 25625.25
$ws.Range("I71").Value = 25625.25
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 230627.25
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -226571.25
$ws.Range("N71").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 29000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 29000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 29000
$ws.Range("N63").Value = -30372

$ws.Range("H66").Value = 29000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 29000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 87000
$ws.Range("N66").Value = -93864

$ws.Range("H102").Value = 3053.25
$ws.Range("I102").Value = 3202.4
$ws.Range("J102").Value = 2804.6667
$ws.Range("K102").Value = 3202.4
$ws.Range("L102").Value = 2804.6667
$ws.Range("M102").Value = -1580.4
$ws.Range("N102").Value = -6048.6667

$ws.Range("H122").Value = 446204.2
$ws.Range("I122").Value = 695686.2
$ws.Range("J122").Value = 2680.6667
$ws.Range("K122").Value = 2087058.6
$ws.Range("L122").Value = 8042.000100000001
$ws.Range("M122").Value = -2084608.6
$ws.Range("N122").Value = -12942.0001

$ws.Range("H132").Value = 3830.8
$ws.Range("I132").Value = 3937.3333
$ws.Range("J132").Value = 3404.6667
$ws.Range("K132").Value = 11811.9999
$ws.Range("L132").Value = 10214.0001
$ws.Range("M132").Value = -9281.999899999999
$ws.Range("N132").Value = -15274.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3409.111
$ws.Range("I132").Value = 2130.5
$ws.Range("J132").Value = 5966.3335
$ws.Range("K132").Value = 6391.5
$ws.Range("L132").Value = 17899.0005
$ws.Range("M132").Value = -3861.5
$ws.Range("N132").Value = -22959.0005

$ws.Range("H136").Value = 4957.275
$ws.Range("I136").Value = 2922.6
$ws.Range("J136").Value = 19200
$ws.Range("K136").Value = 8767.799999999999
$ws.Range("L136").Value = 57600
$ws.Range("M136").Value = -6217.799999999999
$ws.Range("N136").Value = -62700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 162207.25
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 162207.25
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 162207.25
$ws.Range("N46").Value = -162669.25

$ws.Range("H134").Value = 162207.25
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 162207.25
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 486621.75
$ws.Range("N134").Value = -491691.75
